# debug_lines.xlsx : add a header/comment row at the top of Sheet1 and a new
# "F" helper column that tells you which source-file line a given cumulative
# line-number falls on, plus a note about the INVALID placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank row above the existing table -------------------
# (everything that used to be on rows 1-8 shifts down to rows 2-9)
$ws.Rows.Item(1).Insert()

# --- 2. Explanatory note above the table, in column E ----------------------
$ws.Range("E1").Value = "AICI pui linia, in functie de intervalul (coloana A) unde se gaseste, in F iti apare linia unde sa te uiti in fisierul respectiv"

# --- 3. Updated line-count values in column D (recomputed source files) ----
$ws.Range("D2").Value = 1101
$ws.Range("D5").Value = 447
$ws.Range("D6").Value = 783
$ws.Range("D7").Value = 596
$ws.Range("D8").Value = 496
$ws.Range("D9").Value = 1

# --- 4. New helper column F: find the line inside the matching interval ----
$ws.Range("F2:F7").Formula = '=IF(IF(E2>0,E2-A2-5,0)<D2,IF(E2>0,E2-A2-5,0),"INVALID")'
$ws.Range("F8").Formula = '=IF(IF(E8>0,E8-A8-5,0)<D8,IF(E8>0,E8-A8-5,0),"INVALID")'
$ws.Range("F9").Formula = '=IF(IF(E9>0,E9-A9-5,0)<D9,IF(E9>0,E9-A9-5,0),"INVALID")'

# --- 5. Trailing remark next to the last row, in column I ------------------
$ws.Range("I9").Value = "1 am pus ca sa nu scrie INVALID la F, dar ar trebui scris nr de linii din <>JS script de QUX"

# --- 6. Make columns A and B a bit wider so labels/comments are readable ---
$ws.Columns.Item(1).ColumnWidth = 8.25
$ws.Columns.Item(2).ColumnWidth = 17.6

# --- 7. Leave the selection where the author left it (on the new F9 cell) --
$ws.Range("F9").Select() | Out-Null
